$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-03-03"

# Update the header label in column I (shared string "2022 (through 03-02)" -> "...03-03")
$ws.Range("I1").Value = "2022 (through 03-03)"

# Update the March data point for 2022 (I4) and the Total row (I14)
$ws.Range("I4").Value = 16
$ws.Range("I14").Value = 316
